$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph text "PLOG0012 – Contratação de fornecimento de
# material nacionalizado" becomes "PLOG0012 – Ressuprimento de material
# nacionalizado", ending up split across three runs (same rPr):
#   "PLOG0012 – " | "Ressuprimento d" | "e material nacionalizado"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "PLOG0012 – Contratação de fornecimento de material nacionalizado",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $base = $rng.Start

    # Replace the part after "PLOG0012 – " (11 chars) with the new wording.
    $rTail = $d.Range($base + 11, $base + 64)
    $rTail.Text = "Ressuprimento de material nacionalizado"

    # Force the middle chunk ("Ressuprimento d") onto its own run by
    # toggling direct character formatting, which splits it away from its
    # neighbours; then clear the formatting again so the three runs end up
    # with identical (original) run properties.
    $rMiddle = $d.Range($base + 11, $base + 11 + 15)
    $rMiddle.Font.Bold = 1
    $rMiddle2 = $d.Range($base + 11, $base + 11 + 15)
    $rMiddle2.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 2: cached page-number field result in the header goes from "3" to
# "2" (the PAGE field's cached/displayed digit). Writing straight to
# Field.Result.Text lands in the wrong XML node in this host, so the cached
# digit is located via the header's Characters collection (which reports the
# true cached field text, unlike Range.Text) and overwritten character by
# character instead.
# ---------------------------------------------------------------------------
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
$hdrRng = $hdr.Range

$count = $hdrRng.Characters.Count
$joined = ""
for ($i = 1; $i -le $count; $i++) {
    $joined = $joined + $hdrRng.Characters.Item($i).Text
}

foreach ($f in $hdrRng.Fields) {
    $old = $f.Result.Text
    if ($old -eq "3") {
        $new = "2"
        $pos = $joined.LastIndexOf($old)
        if ($pos -ge 0) {
            for ($k = 0; $k -lt $old.Length; $k++) {
                $charIndex = $pos + $k + 1   # Characters is 1-based
                $hdrRng.Characters.Item($charIndex).Text = $new.Substring($k, 1)
            }
        }
    }
}
